$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D/E updates for simple rows ---
$ws.Range("D2").Value = "27.853.02"
$ws.Range("E2").Value = "  +2.60%  "
$ws.Range("D3").Value = "1.668.03"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.59"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.61%  "
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0621"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0879"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.30%  "
$ws.Range("D12").Value = "1.903.48"
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("D13").Value = "1.665.58"
$ws.Range("E13").Value = "  -0.77%  "
$ws.Range("E14").Value = "  -1.33%  "
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "251.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.22%  "
$ws.Range("D18").Value = "27.840.67"
$ws.Range("E18").Value = "  +2.66%  "
$ws.Range("D19").Value = "0.0₃0731"
$ws.Range("E19").Value = "  -1.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.95%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("E22").Value = "  -1.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.60%  "
$ws.Range("E26").Value = "  -2.98%  "
$ws.Range("E27").Value = "  -0.62%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  -0.74%  "
$ws.Range("E30").Value = "  +5.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0500"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("E32").Value = "  -0.33%  "
$ws.Range("E33").Value = "  -2.85%  "
$ws.Range("D34").Value = "1.418.19"
$ws.Range("E34").Value = "  -8.03%  "
$ws.Range("E35").Value = "  -5.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.39"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.580"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.60%  "
$ws.Range("E39").Value = "  -1.41%  "
$ws.Range("E40").Value = "  -3.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.51"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("E43").Value = "  -1.17%  "
$ws.Range("D44").Value = "1.811.48"
$ws.Range("E44").Value = "  -0.87%  "
$ws.Range("E47").Value = "  +4.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.18%  "
$ws.Range("E49").Value = "  -0.55%  "

# --- Rows with Coin/Link swaps (45, 46, 51) ---
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.791"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.33%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.17%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.51%  "
